$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1511.8182
$ws.Range("I6").Value = 312.75
$ws.Range("J6").Value = 2197
$ws.Range("K6").Value = 938.25
$ws.Range("L6").Value = 6591
$ws.Range("M6").Value = -826.25
$ws.Range("N6").Value = -6815

$ws.Range("H33").Value = 736.13043
$ws.Range("I33").Value = 634.375
$ws.Range("K33").Value = 634.375
$ws.Range("M33").Value = -405.375

$ws.Range("H101").Value = 5727.2
$ws.Range("I101").Value = 318
$ws.Range("K101").Value = 954
$ws.Range("M101").Value = 668

$ws.Range("H132").Value = 20008760
$ws.Range("I132").Value = 22230846
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 66692538
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -66690008
$ws.Range("N132").Value = -35060

$ws.Range("H138").Value = 5212.675
$ws.Range("I138").Value = 2543.32
$ws.Range("J138").Value = 6426.018
$ws.Range("K138").Value = 7629.960000000001
$ws.Range("L138").Value = 19278.054
$ws.Range("M138").Value = -2489.960000000001
$ws.Range("N138").Value = -29558.054

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10420.19
$ws.Range("I32").Value = 9183.779
$ws.Range("J32").Value = 18015.285
$ws.Range("K32").Value = 9183.779
$ws.Range("L32").Value = 18015.285
$ws.Range("M32").Value = -8896.779
$ws.Range("N32").Value = -18589.285

$ws.Range("H74").Value = 2575.5715
$ws.Range("I74").Value = 1839.8
$ws.Range("J74").Value = 4415
$ws.Range("K74").Value = 1839.8
$ws.Range("L74").Value = 4415
$ws.Range("M74").Value = -965.8
$ws.Range("N74").Value = -6163

$ws.Range("H77").Value = 2575.5715
$ws.Range("I77").Value = 1839.8
$ws.Range("J77").Value = 4415
$ws.Range("K77").Value = 9199
$ws.Range("L77").Value = 22075
$ws.Range("M77").Value = -4831
$ws.Range("N77").Value = -30811

$ws.Range("H122").Value = 3247.2856
$ws.Range("I122").Value = 2142.3572
$ws.Range("K122").Value = 6427.071599999999
$ws.Range("M122").Value = -3977.071599999999

$ws.Range("H132").Value = 1827.5676
$ws.Range("I132").Value = 1646.3429
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 4939.028700000001
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -2409.028700000001
$ws.Range("N132").Value = -20057

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3083.6978
$ws.Range("I134").Value = 2994.8462
$ws.Range("J134").Value = 3950
$ws.Range("K134").Value = 8984.5386
$ws.Range("L134").Value = 11850
$ws.Range("M134").Value = -6449.5386
$ws.Range("N134").Value = -16920

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2867.6626
$ws.Range("I31").Value = 2036.2273
$ws.Range("J31").Value = 3883.861
$ws.Range("K31").Value = 2036.2273
$ws.Range("L31").Value = 3883.861
$ws.Range("M31").Value = -1741.2273
$ws.Range("N31").Value = -4473.861

$ws.Range("H34").Value = 2867.6626
$ws.Range("I34").Value = 2036.2273
$ws.Range("J34").Value = 3883.861
$ws.Range("K34").Value = 2036.2273
$ws.Range("L34").Value = 3883.861
$ws.Range("M34").Value = -1834.2273
$ws.Range("N34").Value = -4287.861

$ws.Range("H122").Value = 2703.875
$ws.Range("I122").Value = 2265.2778
$ws.Range("K122").Value = 6795.8334
$ws.Range("M122").Value = -4345.8334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 5659.0835
$ws.Range("I56").Value = 5659.0835
$ws.Range("K56").Value = 5659.0835
$ws.Range("M56").Value = -5129.0835

$ws.Range("H122").Value = 1290.5
$ws.Range("I122").Value = 362.57144
$ws.Range("J122").Value = 2589.6
$ws.Range("K122").Value = 3263.14296
$ws.Range("L122").Value = 23306.4
$ws.Range("M122").Value = -813.1429600000001
$ws.Range("N122").Value = -28206.4

$ws.Range("H130").Value = 2169
$ws.Range("I130").Value = 1538
$ws.Range("K130").Value = 4614
$ws.Range("M130").Value = 406

$ws.Range("H131").Value = 1354.7878
$ws.Range("J131").Value = 1135.2858
$ws.Range("L131").Value = 3405.8574
$ws.Range("N131").Value = -13485.8574

$ws.Range("H133").Value = 7384.4443
$ws.Range("I133").Value = 11615
$ws.Range("J133").Value = 4000
$ws.Range("K133").Value = 34845
$ws.Range("L133").Value = 12000
$ws.Range("M133").Value = -29785
$ws.Range("N133").Value = -22120

$ws.Range("H134").Value = 4921.375
$ws.Range("I134").Value = 3753.3333
$ws.Range("J134").Value = 5622.2
$ws.Range("K134").Value = 11259.9999
$ws.Range("L134").Value = 16866.6
$ws.Range("M134").Value = -6189.999899999999
$ws.Range("N134").Value = -27006.6

$ws.Range("H139").Value = 35720820
$ws.Range("I139").Value = 35720820
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 107162460
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -107157320
$ws.Range("N139").ClearContents()

$ws.Range("H140").Value = 20837824
$ws.Range("I140").Value = 33334520
$ws.Range("J140").Value = 9998
$ws.Range("K140").Value = 100003560
$ws.Range("L140").Value = 29994
$ws.Range("M140").Value = -99998380
$ws.Range("N140").Value = -40354

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5013.923
$ws.Range("I122").Value = 3873.8064
$ws.Range("J122").Value = 6696.952
$ws.Range("K122").Value = 11621.4192
$ws.Range("L122").Value = 20090.856
$ws.Range("M122").Value = -9171.4192
$ws.Range("N122").Value = -24990.856

$ws.Range("H126").Value = 913410.4
$ws.Range("J126").Value = 1115557.1
$ws.Range("L126").Value = 3346671.3
$ws.Range("N126").Value = -3351611.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3236.0715
$ws.Range("I122").Value = 2547.3845
$ws.Range("J122").Value = 3832.9333
$ws.Range("K122").Value = 7642.1535
$ws.Range("L122").Value = 11498.7999
$ws.Range("M122").Value = -5192.1535
$ws.Range("N122").Value = -16398.7999

$ws.Range("H132").Value = 2232.6123
$ws.Range("I132").Value = 1639.9429
$ws.Range("J132").Value = 3714.2856
$ws.Range("K132").Value = 4919.8287
$ws.Range("L132").Value = 11142.8568
$ws.Range("M132").Value = -2389.8287
$ws.Range("N132").Value = -16202.8568

$ws.Range("H136").Value = 2691.6365
$ws.Range("I136").Value = 1814.7675
$ws.Range("J136").Value = 5833.75
$ws.Range("K136").Value = 5444.3025
$ws.Range("L136").Value = 17501.25
$ws.Range("M136").Value = -2894.3025
$ws.Range("N136").Value = -22601.25

$ws.Range("H140").Value = 29642.857
$ws.Range("J140").Value = 29642.857
$ws.Range("L140").Value = 29642.857
$ws.Range("N140").Value = -40002.857

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 964.4167
$ws.Range("I113").Value = 121.666664
$ws.Range("J113").Value = 1807.1666
$ws.Range("K113").Value = 364.999992
$ws.Range("L113").Value = 5421.4998
$ws.Range("M113").Value = 1805.000008
$ws.Range("N113").Value = -9761.4998

$ws.Range("H122").Value = 1967.6666
$ws.Range("J122").Value = 2777.4
$ws.Range("L122").Value = 8332.200000000001
$ws.Range("N122").Value = -13232.2

$ws.Range("H136").Value = 2188.9556
$ws.Range("I136").Value = 1751.3
$ws.Range("K136").Value = 5253.9
$ws.Range("M136").Value = -2703.9

$ws.Range("H141").Value = 29666.666
$ws.Range("J141").Value = 29666.666
$ws.Range("L141").Value = 29666.666
$ws.Range("N141").Value = -40026.666
